$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation for the Price/Volume columns so numeric-looking
# strings (e.g. "1.00", "35.19") are not auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '42.255.59'
$ws.Range('E2').Value = '  -0.99%  '
$ws.Range('D3').Value = '2.268.02'
$ws.Range('E3').Value = '  -1.40%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '306.23'
$ws.Range('E5').Value = '  -0.74%  '
$ws.Range('D6').Value = '97.27'
$ws.Range('E6').Value = '  +0.62%  '
$ws.Range('E7').Value = '  -0.99%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  -1.15%  '
$ws.Range('D10').Value = '35.19'
$ws.Range('E10').Value = '  -1.41%  '
$ws.Range('E11').Value = '  -2.44%  '
$ws.Range('E12').Value = '  -0.11%  '
$ws.Range('D13').Value = '6.99'
$ws.Range('E13').Value = '  +3.12%  '
$ws.Range('D14').Value = '2.620.31'
$ws.Range('E14').Value = '  -1.49%  '
$ws.Range('D15').Value = '14.75'
$ws.Range('E15').Value = '  +0.91%  '
$ws.Range('D16').Value = '2.272.41'
$ws.Range('E16').Value = '  -1.67%  '
$ws.Range('D17').Value = '0.794'
$ws.Range('E17').Value = '  -0.44%  '
$ws.Range('D18').Value = '42.124.37'
$ws.Range('E18').Value = '  -1.01%  '
$ws.Range('D19').Value = '12.34'
$ws.Range('E19').Value = '  -3.95%  '
$ws.Range('D20').Value = '0.0₃0907'
$ws.Range('E20').Value = '  -1.68%  '
$ws.Range('E21').Value = '  -0.56%  '
$ws.Range('D22').Value = '67.88'
$ws.Range('E22').Value = '  -0.61%  '
$ws.Range('D23').Value = '237.98'
$ws.Range('E23').Value = '  -2.87%  '
$ws.Range('B24').Value = 'ImmutableX'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D24').Value = '1.97'
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').Value = '2.57'
$ws.Range('E25').Value = '  -1.72%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('E27').Value = '  -3.21%  '
$ws.Range('D28').Value = '38.17'
$ws.Range('E28').Value = '  +2.35%  '
$ws.Range('D29').Value = '9.58'
$ws.Range('E29').Value = '  -1.73%  '
$ws.Range('E30').Value = '  +0.31%  '
$ws.Range('D31').Value = '163.05'
$ws.Range('E31').Value = '  +0.85%  '
$ws.Range('D32').Value = '5.26'
$ws.Range('E32').Value = '  -2.67%  '
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('D34').Value = '3.17'
$ws.Range('E34').Value = '  +1.56%  '
$ws.Range('D35').Value = '17.77'
$ws.Range('E35').Value = '  +1.92%  '
$ws.Range('D36').Value = '0.0739'
$ws.Range('E36').Value = '  -2.62%  '
$ws.Range('E37').Value = '  -0.82%  '
$ws.Range('E38').Value = '  -4.51%  '
$ws.Range('D39').Value = '1.83'
$ws.Range('E39').Value = '  -1.44%  '
$ws.Range('E40').Value = '  -1.69%  '
$ws.Range('E41').Value = '  -4.11%  '
$ws.Range('E42').Value = '  +2.81%  '
$ws.Range('D43').Value = '19.14'
$ws.Range('E43').Value = '  -3.93%  '
$ws.Range('D44').Value = '1.955.60'
$ws.Range('E44').Value = '  -3.25%  '
$ws.Range('E45').Value = '  -1.62%  '
$ws.Range('E46').Value = '  -3.25%  '
$ws.Range('D47').Value = '9.87'
$ws.Range('E47').Value = '  -4.25%  '
$ws.Range('D48').Value = '53.95'
$ws.Range('E48').Value = '  +0.01%  '
$ws.Range('D49').Value = '92.51'
$ws.Range('E49').Value = '  -0.88%  '
$ws.Range('D50').Value = '71.75'
$ws.Range('E50').Value = '  -2.14%  '
$ws.Range('E51').Value = '  -2.40%  '

# Restore default (General) formatting so no stray cell styles are introduced.
$ws.Range("D2:E51").ClearFormats()
